$wb = $excel.ActiveWorkbook
$wsHist = $wb.Worksheets.Item("Historical")
$wsRT   = $wb.Worksheets.Item("Real Time")

# --- Value edits -----------------------------------------------------

# "Historical" sheet: the account number text "1619205" (shared string,
# used by rows 20-25 in column C) is corrected to "1619215".
$wsHist.Range("C20").Value = "1619215"
$wsHist.Range("C21").Value = "1619215"
$wsHist.Range("C22").Value = "1619215"
$wsHist.Range("C23").Value = "1619215"
$wsHist.Range("C24").Value = "1619215"
$wsHist.Range("C25").Value = "1619215"

# "Real Time" sheet: row 4 GL/Difference balances update.
$wsRT.Range("H4").Value = 60000
$wsRT.Range("J4").Value = 50000

# "Real Time" sheet: row 5 account number corrected from 1619205 to 1619215.
$wsRT.Range("C5").Value = 1619215

# --- Selection / view state -------------------------------------------
# Historical sheet ends up scrolled to A1 with C25 selected.
[void]$wsHist.Range("C25").Select()

# Real Time sheet stays the active tab, scrolled to A1, with J4 selected.
[void]$wsRT.Range("J4").Select()
